$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D11"  = -7.852499999999989
    "B12"  = 4.829399999999998
    "D23"  = -8.003
    "D28"  = -8.357999999999997
    "B32"  = 6.4575
    "D32"  = -8.511899999999986
    "D34"  = -8.096500000000002
    "B36"  = 9.188700000000004
    "B38"  = 4.864399999999998
    "D42"  = -9.101799999999992
    "B46"  = 5.769200000000003
    "B54"  = 4.882099999999999
    "D54"  = -8.350799999999994
    "B55"  = 6.046999999999998
    "B67"  = 5.309599999999993
    "B69"  = 5.304699999999996
    "B72"  = 5.196600000000003
    "B91"  = 5.149399999999998
    "D97"  = -8.48379999999999
    "B99"  = 4.626800000000004
    "D99"  = -8.644100000000003
    "D101" = -7.785999999999999
    "B104" = 9.832400000000003
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
